$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.745631333333333
$ws.Range("H2").Value = 11.236894
$ws.Range("I2").Value = 0.1419671142338921
$ws.Range("J2").Value = 0.1419671142338921
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 78.96470183336154
$ws.Range("R2").Value = 710.6823165002539
$ws.Range("S2").Value = 0.008114925637129141
$ws.Range("T2").Value = 0.008114925637129141
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.745631333333333
$ws.Range("H3").Value = 11.236894
$ws.Range("I3").Value = 0.1419671142338921
$ws.Range("J3").Value = 0.1419671142338921
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 1129.682966983858
$ws.Range("R3").Value = 10167.14670285472
$ws.Range("S3").Value = 0.1160935589923593
$ws.Range("T3").Value = 0.1160935589923593
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.745631333333333
$ws.Range("H4").Value = 11.236894
$ws.Range("I4").Value = 0.1419671142338921
$ws.Range("J4").Value = 0.1419671142338921
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 172.8056367226229
$ws.Range("R4").Value = 1555.250730503606
$ws.Range("S4").Value = 0.0177586296044036
$ws.Range("T4").Value = 0.0177586296044036
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.68955666666667
$ws.Range("H5").Value = 47.06867
$ws.Range("I5").Value = 0.5946663954227359
$ws.Range("J5").Value = 0.5946663954227359
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 330.7643101592744
$ws.Range("R5").Value = 2976.87879143347
$ws.Range("S5").Value = 0.03399148883032725
$ws.Range("T5").Value = 0.03399148883032725
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.68955666666667
$ws.Range("H6").Value = 47.06867
$ws.Range("I6").Value = 0.5946663954227359
$ws.Range("J6").Value = 0.5946663954227359
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 4731.972623180757
$ws.Range("R6").Value = 42587.75360862682
$ws.Range("S6").Value = 0.4862882409798379
$ws.Range("T6").Value = 0.4862882409798379
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.68955666666667
$ws.Range("H7").Value = 47.06867
$ws.Range("I7").Value = 0.5946663954227359
$ws.Range("J7").Value = 0.5946663954227359
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 723.8416139759812
$ws.Range("R7").Value = 6514.574525783831
$ws.Range("S7").Value = 0.07438666561257085
$ws.Range("T7").Value = 0.07438666561257085
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.948607666666667
$ws.Range("H8").Value = 20.845823
$ws.Range("I8").Value = 0.263366490343372
$ws.Range("J8").Value = 0.263366490343372
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 146.4892520714381
$ws.Range("R8").Value = 1318.403268642943
$ws.Range("S8").Value = 0.01505418699239811
$ws.Range("T8").Value = 0.01505418699239811
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.948607666666667
$ws.Range("H9").Value = 20.845823
$ws.Range("I9").Value = 0.263366490343372
$ws.Range("J9").Value = 0.263366490343372
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 2095.701105293007
$ws.Range("R9").Value = 18861.30994763706
$ws.Range("S9").Value = 0.215367857184982
$ws.Range("T9").Value = 0.215367857184982
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.948607666666667
$ws.Range("H10").Value = 20.845823
$ws.Range("I10").Value = 0.263366490343372
$ws.Range("J10").Value = 0.263366490343372
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 320.5757495373808
$ws.Range("R10").Value = 2885.181745836427
$ws.Range("S10").Value = 0.03294444616599191
$ws.Range("T10").Value = 0.03294444616599191
